$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Price/Volume columns to text so Excel does not coerce numeric-looking
# strings (e.g. "1.00", "0.208") into numbers and strip formatting/precision.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

$ws.Range("D2").Value = "69.782.49"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "3.548.98"
$ws.Range("E3").Value = "  -1.90%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "591.16"
$ws.Range("E5").Value = "  -2.01%  "
$ws.Range("D6").Value = "196.17"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -2.32%  "
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("D9").Value = "0.208"
$ws.Range("E9").Value = "  +0.89%  "
$ws.Range("D10").Value = "0.626"
$ws.Range("E10").Value = "  -3.12%  "
$ws.Range("D11").Value = "52.83"
$ws.Range("E11").Value = "  -1.08%  "
$ws.Range("E12").Value = "  -3.79%  "
$ws.Range("D13").Value = "9.32"
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").Value = "4.120.98"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "644.71"
$ws.Range("E15").Value = "  +7.59%  "
$ws.Range("D16").Value = "69.820.54"
$ws.Range("E16").Value = "  -1.02%  "
$ws.Range("D17").Value = "12.57"
$ws.Range("E17").Value = "  -3.38%  "
$ws.Range("D18").Value = "3.553.31"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("D20").Value = "18.49"
$ws.Range("E20").Value = "  -2.82%  "
$ws.Range("D21").Value = "0.965"
$ws.Range("E21").Value = "  -3.05%  "
$ws.Range("D22").Value = "18.11"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  +6.05%  "
$ws.Range("D24").Value = "103.28"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("E25").Value = "  -4.18%  "
$ws.Range("D26").Value = "2.94"
$ws.Range("E26").Value = "  -2.29%  "
$ws.Range("D27").Value = "10.27"
$ws.Range("E27").Value = "  -4.26%  "
$ws.Range("D28").Value = "9.62"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").Value = "33.28"
$ws.Range("E29").Value = "  -1.57%  "
$ws.Range("D30").Value = "4.28"
$ws.Range("E30").Value = "  -7.58%  "
$ws.Range("D31").Value = "6.84"
$ws.Range("E31").Value = "  -5.31%  "
$ws.Range("D32").Value = "11.79"
$ws.Range("E32").Value = "  -3.82%  "
$ws.Range("D33").Value = "0.111"
$ws.Range("E33").Value = "  -6.44%  "
$ws.Range("D34").Value = "61.77"
$ws.Range("E34").Value = "  -2.54%  "
$ws.Range("D35").Value = "3.736.39"
$ws.Range("E35").Value = "  -4.42%  "
$ws.Range("D36").Value = "0.0₃0815"
$ws.Range("E36").Value = "  -8.03%  "
$ws.Range("D37").Value = "0.999"
$ws.Range("E37").Value = "  -0.22%  "
$ws.Range("D38").Value = "3.67"
$ws.Range("E38").Value = "  +3.97%  "
$ws.Range("D39").Value = "515.10"
$ws.Range("E39").Value = "  -3.90%  "
$ws.Range("E40").Value = "  -4.12%  "
$ws.Range("D41").Value = "0.373"
$ws.Range("E41").Value = "  -4.23%  "
$ws.Range("D42").Value = "35.13"
$ws.Range("E42").Value = "  -4.81%  "
$ws.Range("D43").Value = "0.134"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("D44").Value = "0.0453"
$ws.Range("E44").Value = "  -2.27%  "
$ws.Range("D45").Value = "3.43"
$ws.Range("E45").Value = "  +1.06%  "
$ws.Range("E46").Value = "  +0.57%  "
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").Value = "1.00"
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "8.29"
$ws.Range("E49").Value = "  -3.38%  "
$ws.Range("D50").Value = "1.77"
$ws.Range("E50").Value = "  +18.75%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").Value = "2.72"
$ws.Range("E51").Value = "  +61.59%  "

# Restore the default style so no stray number-format style attribute
# is left behind on these cells.
$dataRange.Style = "Normal"

